$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look fully numeric must be forced to Text format
# first, otherwise Excel auto-converts the assigned string into a real number
# (stripping formatting such as leading/trailing zeros).
$textFormatCells = @{
    "D5" = "271.44"
    "D6" = "92.42"
    "D10" = "45.78"
    "D11" = "0.0970"
    "D12" = "8.32"
    "D19" = "0.0000105"
    "D20" = "6.12"
    "D21" = "70.86"
    "D22" = "2.34"
    "D23" = "234.95"
    "D24" = "9.20"
    "D26" = "11.47"
    "D27" = "2.51"
    "D28" = "3.55"
    "D29" = "41.40"
    "D31" = "172.76"
    "D32" = "0.0919"
    "D33" = "21.01"
    "D34" = "5.51"
    "D37" = "0.0351"
    "D39" = "3.54"
    "D42" = "2.18"
    "D43" = "63.91"
    "D45" = "0.100"
    "D46" = "8.37"
    "D47" = "100.37"
    "D48" = "1.15"
    "D49" = "1.19"
}

foreach ($ref in $textFormatCells.Keys) {
    $ws.Range($ref).NumberFormat = "@"
}
foreach ($ref in $textFormatCells.Keys) {
    $ws.Range($ref).Value = $textFormatCells[$ref]
}
foreach ($ref in $textFormatCells.Keys) {
    $ws.Range($ref).Style = "Normal"
}

# Cells whose new values contain multiple "." separators are never
# auto-converted to numbers by Excel, so they can be assigned directly.
$directCells = @{
    "D2" = "43.958.19"
    "D3" = "2.245.38"
    "D14" = "2.582.64"
    "D16" = "2.247.88"
    "D18" = "43.902.22"
    "D51" = "2.468.43"
}
foreach ($ref in $directCells.Keys) {
    $ws.Range($ref).Value = $directCells[$ref]
}

# Volume(1h) percentage text values (kept as text with the original padding spaces)
$volumeCells = @{
    "E2" = "  +0.93%  "
    "E3" = "  +2.79%  "
    "E4" = "  +0.07%  "
    "E5" = "  +5.05%  "
    "E6" = "  +13.87%  "
    "E7" = "  +1.20%  "
    "E8" = "  -0.03%  "
    "E9" = "  +6.56%  "
    "E10" = "  +6.60%  "
    "E11" = "  +5.81%  "
    "E12" = "  +20.01%  "
    "E13" = "  +1.86%  "
    "E14" = "  +2.93%  "
    "E15" = "  +6.71%  "
    "E16" = "  +3.56%  "
    "E17" = "  +5.14%  "
    "E18" = "  +1.01%  "
    "E19" = "  +3.22%  "
    "E20" = "  +3.94%  "
    "E21" = "  +1.50%  "
    "E22" = "  -1.19%  "
    "E23" = "  +2.24%  "
    "E24" = "  +4.59%  "
    "E25" = "  -0.16%  "
    "E26" = "  +7.96%  "
    "E27" = "  +13.36%  "
    "E28" = "  +5.60%  "
    "E29" = "  -2.52%  "
    "E30" = "  +0.79%  "
    "E31" = "  -0.10%  "
    "E32" = "  +5.83%  "
    "E33" = "  +3.32%  "
    "E34" = "  +4.58%  "
    "E35" = "  +2.04%  "
    "E36" = "  +0.92%  "
    "E37" = "  +0.28%  "
    "E38" = "  -2.90%  "
    "E39" = "  +25.31%  "
    "E40" = "  +0.97%  "
    "E41" = "  +14.69%  "
    "E42" = "  +4.61%  "
    "E43" = "  +2.21%  "
    "E44" = "  -0.25%  "
    "E45" = "  +2.24%  "
    "E46" = "  +2.69%  "
    "E47" = "  -0.11%  "
    "E48" = "  +4.53%  "
    "E49" = "  +1.78%  "
    "E50" = "  +1.07%  "
    "E51" = "  +3.05%  "
}
foreach ($ref in $volumeCells.Keys) {
    $ws.Range($ref).Value = $volumeCells[$ref]
}
